# Regenerate the "K" column (G) of save_data to hold actual strikeout
# counts (K) instead of the previous Strike# totals.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 5
    3  = 8
    4  = 7
    5  = 4
    6  = 4
    7  = 10
    8  = 5
    9  = 9
    10 = 6
    11 = 4
    12 = 6
    13 = 5
    14 = 1
    15 = 6
    16 = 6
    17 = 3
    18 = 2
    19 = 8
    20 = 1
    21 = 8
    22 = 6
    23 = 4
    24 = 1
    25 = 9
    26 = 3
    27 = 10
    28 = 5
    29 = 7
    30 = 4
    31 = 7
    32 = 6
    33 = 2
    34 = 6
    35 = 2
    36 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
